$d = $word.ActiveDocument

# --- 1. Remove the "Translated with www.DeepL.com/Translator (free version)"
#        paragraph and merge its (now empty) paragraph with the preceding
#        empty paragraph, leaving a single empty paragraph that still holds
#        the _GoBack bookmark. ---
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
$last.Range.Find.Execute("Translated with www.DeepL.com/Translator (free version)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$count = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($count - 1)
$r = $secondLast.Range
$r.MoveEnd(1, 1)
$r.Delete()

# --- 2. Mark the built-in "Default Paragraph Font" and "Normal Table"
#        styles as QuickStyle (w:qFormat) styles. ---
$d.Styles.Item("Default Paragraph Font").QuickStyle = $true
$d.Styles.Item("Normal Table").QuickStyle = $true
